$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Drop the stray _GoBack bookmark that currently sits right after
#    the "DP" run (it marks wherever the last edit in the source doc
#    happened to be, which is no longer where we are about to type).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Replace "Recursos Humanos " with "Pessoal " inside the sentence
#    "Foi encaminhada à Diretoria de Recursos Humanos do CBMAM, ..."
# ------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Recursos Humanos ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$rng.Text = "Pessoal "

# Nudge formatting on/off so the newly typed text lands in its own
# run instead of being re-merged into the preceding run (mirrors how
# Word keeps separately-typed text in its own <w:r>).
$rng.Bold = 1
$rng.Bold = 0

# ------------------------------------------------------------------
# 3) Word drops a fresh _GoBack bookmark exactly where the user just
#    finished typing -- right after "Pessoal " and before "do CBMAM".
# ------------------------------------------------------------------
$newSpot = $d.Range($rng.End, $rng.End)
$d.Bookmarks.Add("_GoBack", $newSpot)
